$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 47, shifting rows 47:137 down to 48:138
$ws.Rows.Item(47).Insert()

# Populate the new row 47 with its data
$ws.Cells.Item(47, 1).Value = 9
$ws.Cells.Item(47, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(47, 3).Value = "Metropolitana"
$ws.Cells.Item(47, 4).Value = 45203
$ws.Cells.Item(47, 5).Value = 13
$ws.Cells.Item(47, 6).Value = 100112005
$ws.Cells.Item(47, 7).Value = "Puerro"
$ws.Cells.Item(47, 8).Value = "Sin especificar"
$ws.Cells.Item(47, 9).Value = "Primera"
$ws.Cells.Item(47, 10).Value = 70
$ws.Cells.Item(47, 11).Value = 8000
$ws.Cells.Item(47, 12).Value = 8000
$ws.Cells.Item(47, 13).Value = 8000
$ws.Cells.Item(47, 14).Value = "`$/paquete 20 unidades"
$ws.Cells.Item(47, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(47, 16).Value = 400
$ws.Cells.Item(47, 17).Value = 20
$ws.Cells.Item(47, 18).Value = "Hortaliza"
